$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (style) from H1 to I1:J1, matching existing header style
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Set data values for I2:J75
$data = New-Object 'object[,]' 74,2
$data[0,0] = 5
$data[0,1] = 6
$data[1,0] = 9
$data[1,1] = 9
$data[2,0] = 7
$data[2,1] = 8
$data[3,0] = 9
$data[3,1] = 9
$data[4,0] = 8
$data[4,1] = 9
$data[5,0] = 9
$data[5,1] = 9
$data[6,0] = 9
$data[6,1] = 9
$data[7,0] = 10
$data[7,1] = 10
$data[8,0] = 9
$data[8,1] = 10
$data[9,0] = 9
$data[9,1] = 9
$data[10,0] = 8
$data[10,1] = 9
$data[11,0] = 9
$data[11,1] = 9
$data[12,0] = 9
$data[12,1] = 9
$data[13,0] = 9
$data[13,1] = 9
$data[14,0] = 8
$data[14,1] = 8
$data[15,0] = 8
$data[15,1] = 8
$data[16,0] = 7
$data[16,1] = 7
$data[17,0] = 9
$data[17,1] = 9
$data[18,0] = 9
$data[18,1] = 9
$data[19,0] = 9
$data[19,1] = 9
$data[20,0] = 7
$data[20,1] = 8
$data[21,0] = 8
$data[21,1] = 9
$data[22,0] = 8
$data[22,1] = 9
$data[23,0] = 7
$data[23,1] = 7
$data[24,0] = 7
$data[24,1] = 8
$data[25,0] = 8
$data[25,1] = 8
$data[26,0] = 6
$data[26,1] = 6
$data[27,0] = 8
$data[27,1] = 8
$data[28,0] = 9
$data[28,1] = 9
$data[29,0] = 8
$data[29,1] = 8
$data[30,0] = 8
$data[30,1] = 8
$data[31,0] = 8
$data[31,1] = 8
$data[32,0] = 8
$data[32,1] = 8
$data[33,0] = 7
$data[33,1] = 8
$data[34,0] = 8
$data[34,1] = 8
$data[35,0] = 8
$data[35,1] = 8
$data[36,0] = 8
$data[36,1] = 8
$data[37,0] = 8
$data[37,1] = 8
$data[38,0] = 7
$data[38,1] = 7
$data[39,0] = 8
$data[39,1] = 8
$data[40,0] = 7
$data[40,1] = 8
$data[41,0] = 8
$data[41,1] = 8
$data[42,0] = 4
$data[42,1] = 5
$data[43,0] = 8
$data[43,1] = 8
$data[44,0] = 8
$data[44,1] = 8
$data[45,0] = 7
$data[45,1] = 8
$data[46,0] = 9
$data[46,1] = 10
$data[47,0] = 8
$data[47,1] = 8
$data[48,0] = 7
$data[48,1] = 7
$data[49,0] = 7
$data[49,1] = 7
$data[50,0] = 8
$data[50,1] = 8
$data[51,0] = 8
$data[51,1] = 8
$data[52,0] = 8
$data[52,1] = 9
$data[53,0] = 8
$data[53,1] = 8
$data[54,0] = 9
$data[54,1] = 9
$data[55,0] = 7
$data[55,1] = 7
$data[56,0] = 7
$data[56,1] = 8
$data[57,0] = 8
$data[57,1] = 8
$data[58,0] = 8
$data[58,1] = 8
$data[59,0] = 8
$data[59,1] = 8
$data[60,0] = 8
$data[60,1] = 8
$data[61,0] = 8
$data[61,1] = 8
$data[62,0] = 8
$data[62,1] = 8
$data[63,0] = 7
$data[63,1] = 7
$data[64,0] = 8
$data[64,1] = 8
$data[65,0] = 8
$data[65,1] = 8
$data[66,0] = 8
$data[66,1] = 8
$data[67,0] = 5
$data[67,1] = 5
$data[68,0] = 4
$data[68,1] = 4
$data[69,0] = 6
$data[69,1] = 6
$data[70,0] = 4
$data[70,1] = 5
$data[71,0] = 6
$data[71,1] = 6
$data[72,0] = 3
$data[72,1] = 3
$data[73,0] = 5
$data[73,1] = 5
$ws.Range("I2:J75").Value = $data
